# Apply updated symbol list values (coinranking snapshot refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Sheet, $Row, $Col, $Text) {
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $Text
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws 2 4 '314.43'
Set-TextValue $ws 2 5 '2.92%'

# Row 3
Set-TextValue $ws 3 4 '39.54'
Set-TextValue $ws 3 5 '3.28%'

# Row 4
Set-TextValue $ws 4 4 '5.113'
Set-TextValue $ws 4 5 '0.42%'

# Row 5
Set-TextValue $ws 5 4 '0.08176'
Set-TextValue $ws 5 5 '1.44%'

# Row 6
Set-TextValue $ws 6 4 '2.052'
Set-TextValue $ws 6 5 '6.18%'

# Row 7
Set-TextValue $ws 7 4 '8.247'
Set-TextValue $ws 7 5 '3.85%'

# Row 8
Set-TextValue $ws 8 2 'GateToken'
Set-TextValue $ws 8 3 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 8 4 '4.276'
Set-TextValue $ws 8 5 '2.13%'

# Row 9
Set-TextValue $ws 9 2 'MXToken'
Set-TextValue $ws 9 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 9 4 '0.9338'
Set-TextValue $ws 9 5 '0.46%'

# Row 10
Set-TextValue $ws 10 2 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 10 3 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 10 4 '0.1412'
Set-TextValue $ws 10 5 '-1.80%'

# Row 11
Set-TextValue $ws 11 2 'WazirX'
Set-TextValue $ws 11 3 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 11 4 '0.1989'
Set-TextValue $ws 11 5 '3.46%'

# Row 12
Set-TextValue $ws 12 2 'MandalaExchangeToken'
Set-TextValue $ws 12 3 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 12 4 '0.09085'
Set-TextValue $ws 12 5 '0.90%'

# Row 13
Set-TextValue $ws 13 2 'BitrueCoin'
Set-TextValue $ws 13 3 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 13 4 '0.03533'
Set-TextValue $ws 13 5 '0.32%'

# Row 14
Set-TextValue $ws 14 2 'BitMartToken'
Set-TextValue $ws 14 3 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 14 4 '0.09817'

# Row 15
Set-TextValue $ws 15 2 'BitForexToken'
Set-TextValue $ws 15 3 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 15 4 '0.001398'
Set-TextValue $ws 15 5 '0.57%'

# Row 16
Set-TextValue $ws 16 2 'TigerCash'
Set-TextValue $ws 16 3 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 16 4 '0.006270'
Set-TextValue $ws 16 5 '2.88%'

# Row 17
Set-TextValue $ws 17 2 'LEO'
Set-TextValue $ws 17 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 17 4 '3.659'
Set-TextValue $ws 17 5 '-1.66%'

# Row 18
Set-TextValue $ws 18 4 '3.318'
Set-TextValue $ws 18 5 '-4.03%'

# Row 19
Set-TextValue $ws 19 4 '0.3459'
Set-TextValue $ws 19 5 '-0.09%'

# Row 20
Set-TextValue $ws 20 4 '0.1304'
Set-TextValue $ws 20 5 '-0.57%'

# Row 21
Set-TextValue $ws 21 4 '4.906'
Set-TextValue $ws 21 5 '2.31%'

# Row 22
Set-TextValue $ws 22 4 '0.2452'
Set-TextValue $ws 22 5 '1.88%'

# Row 23
Set-TextValue $ws 23 4 '0.04338'
Set-TextValue $ws 23 5 '-0.24%'

# Row 24
Set-TextValue $ws 24 4 '0.001226'
Set-TextValue $ws 24 5 '-0.35%'

# Row 25
Set-TextValue $ws 25 4 '0.004782'
Set-TextValue $ws 25 5 '16.19%'

# Row 26
Set-TextValue $ws 26 4 '0.0001302'
Set-TextValue $ws 26 5 '0.01%'

# Row 27
Set-TextValue $ws 27 4 '0.0004003'
Set-TextValue $ws 27 5 '-9.99%'

# Row 39
Set-TextValue $ws 39 4 '0.02233'
Set-TextValue $ws 39 5 '8.22%'

# Row 40
Set-TextValue $ws 40 4 '0.05248'
Set-TextValue $ws 40 5 '4.64%'

# Row 41
Set-TextValue $ws 41 4 '0.007573'
Set-TextValue $ws 41 5 '1.23%'

# Row 42
Set-TextValue $ws 42 4 '0.009749'
Set-TextValue $ws 42 5 '-3.65%'

# Row 43
Set-TextValue $ws 43 4 '0.1379'
Set-TextValue $ws 43 5 '2.43%'

# Row 44
Set-TextValue $ws 44 4 '0.002143'
Set-TextValue $ws 44 5 '0.02%'

# Row 45
Set-TextValue $ws 45 4 '0.009788'
Set-TextValue $ws 45 5 '9.72%'

# Row 46
Set-TextValue $ws 46 4 '0.00006414'
Set-TextValue $ws 46 5 '3.58%'

# Row 47
Set-TextValue $ws 47 5 '-0.04%'

# Row 48
Set-TextValue $ws 48 4 '0.002769'
Set-TextValue $ws 48 5 '-1.86%'

# Row 49
Set-TextValue $ws 49 5 '-25.09%'

# Row 50
Set-TextValue $ws 50 4 '0.00002101'
Set-TextValue $ws 50 5 '-0.04%'

# Row 51
Set-TextValue $ws 51 4 '0.0002001'
Set-TextValue $ws 51 5 '-0.04%'
